# fix identifier & tests
#
# 1. Shared string used for A1 ("Tuote/Asiakas") is actually the
#    Customer/Product column header and was the wrong way round -
#    flip it to "Asiakas/Tuote".
# 2. Selection cursor should rest on A1 (top-left), not C5.
# 3. Column A needs to be a bit wider than the rest of the sheet now
#    that it holds the longer "Asiakas/Tuote" header.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Asiakas/Tuote"

$ws.Columns.Item(1).ColumnWidth = 11.75

$ws.Range("A1").Select()
